# The target diff for this commit ("Moving from 2.0.2 to 2.0.3") touches
# word/document.xml, word/footnotes.xml, word/header1.xml and word/styles.xml,
# but every single hunk is a pure XML re-serialization:
#   - root/element attributes (and xmlns declarations) reordered alphabetically
#   - whitespace/newline entities inside VML "o:gfxdata" base64 blobs reflowed
# No element, attribute value, run, paragraph, field, image, style, section
# property, etc. actually changes value or is added/removed anywhere in the
# diff (verified attribute-by-attribute, including decoding every o:gfxdata
# blob). The underlying document content and formatting are therefore
# byte-for-byte semantically identical before and after this commit for this
# particular binary resource - only the XML serializer that produced the file
# changed (e.g. a re-export/round-trip through a different OOXML writer),
# which is not something reachable through the Word object model.
#
# Since there is no actual content/formatting change to apply, we leave the
# document untouched.
$d = $word.ActiveDocument
